$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.098.37"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "1.790.92"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("E6").Value = "  -0.59%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.19"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -1.90%  "
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0692"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0941"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.56%  "
$ws.Range("D12").Value = "2.049.19"
$ws.Range("E12").Value = "  +0.21%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.55"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +5.96%  "
$ws.Range("D14").Value = "1.793.14"
$ws.Range("E14").Value = "  +0.81%  "
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "34.096.12"
$ws.Range("E16").Value = "  +0.15%  "
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.05"
$ws.Range("D18").ClearFormats()
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.55"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -0.42%  "
$ws.Range("D20").Value = "0.0₃0782"
$ws.Range("E20").Value = "  -1.08%  "
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("E24").Value = "  -3.20%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "162.59"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.15%  "
$ws.Range("E26").Value = "  +1.61%  "
$ws.Range("E27").Value = "  -0.21%  "
$ws.Range("E28").Value = "  +1.12%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.24"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("E31").Value = "  +1.15%  "
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("E34").Value = "  +1.01%  "
$ws.Range("D35").Value = "1.415.74"
$ws.Range("E35").Value = "  +1.39%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.642"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.13%  "
$ws.Range("E37").Value = "  +2.42%  "
$ws.Range("E38").Value = "  -1.30%  "
$ws.Range("E39").Value = "  +5.97%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "81.05"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.16%  "
$ws.Range("E41").Value = "  +0.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.922"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +0.63%  "
$ws.Range("E43").Value = "  -0.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.40"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +7.11%  "
$ws.Range("E45").Value = "  +3.22%  "
$ws.Range("E46").Value = "  +2.00%  "
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("E48").Value = "  -5.94%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "107.38"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.76%  "
$ws.Range("D50").Value = "1.949.71"
$ws.Range("E50").Value = "  +0.12%  "
$ws.Range("E51").Value = "  +0.11%  "
